# Updated symbol list on Mon Jan  9 06:33:35 UTC 2023 with GitHub Actions
# Applies the refreshed Price (column D) / Volume(1h) (column E) figures.
# Values are plain-text cells in the source sheet, so each assignment is
# prefixed with a leading apostrophe to force Excel to keep it as text
# (preventing numeric auto-conversion / precision drift and preserving
# formatting such as trailing zeros and the trailing "%" sign).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

Set-Text "D2"  "278.41"
Set-Text "E2"  "6.59%"

Set-Text "D3"  "27.16"
Set-Text "E3"  "0.53%"

Set-Text "D4"  "4.853"
Set-Text "E4"  "3.41%"

Set-Text "D5"  "0.06258"
Set-Text "E5"  "0.65%"

Set-Text "D6"  "6.885"
Set-Text "E6"  "1.93%"

Set-Text "D7"  "0.8794"
Set-Text "E7"  "3.29%"

Set-Text "D8"  "0.9429"
Set-Text "E8"  "2.82%"

Set-Text "D9"  "0.1449"
Set-Text "E9"  "3.41%"

Set-Text "D10" "0.05221"
Set-Text "E10" "7.39%"

Set-Text "D11" "0.07328"
Set-Text "E11" "3.39%"

Set-Text "D12" "0.03163"
Set-Text "E12" "2.00%"

Set-Text "D13" "0.09053"
Set-Text "E13" "-0.02%"

Set-Text "D14" "0.001553"
Set-Text "E14" "1.61%"

Set-Text "D15" "0.0006270"
Set-Text "E15" "1.44%"

Set-Text "D16" "0.005953"
Set-Text "E16" "-3.16%"

Set-Text "E17" "0.30%"

Set-Text "D18" "3.266"
Set-Text "E18" "2.69%"

Set-Text "D19" "2.286"
Set-Text "E19" "5.55%"

Set-Text "E20" "-0.61%"

Set-Text "D21" "0.1312"
Set-Text "E21" "0.07%"

Set-Text "D22" "3.847"
Set-Text "E22" "-6.08%"

Set-Text "D23" "0.04319"
Set-Text "E23" "1.84%"

Set-Text "E24" "-2.10%"

Set-Text "D25" "0.004275"
Set-Text "E25" "4.79%"

Set-Text "D26" "0.0001199"
Set-Text "E26" "-0.12%"

Set-Text "D27" "0.0001691"
Set-Text "E27" "3.09%"

Set-Text "D40" "0.04018"
Set-Text "E40" "0.97%"

Set-Text "D41" "0.006414"
Set-Text "E41" "55.29%"

Set-Text "E42" "3.61%"

Set-Text "D43" "0.002108"
Set-Text "E43" "-4.68%"

Set-Text "D44" "0.01217"
Set-Text "E44" "-12.36%"

Set-Text "D45" "0.00005107"
Set-Text "E45" "-1.08%"

Set-Text "E46" "-0.12%"

Set-Text "D47" "2.369"
Set-Text "E47" "811.53%"

Set-Text "D49" "0.00002099"
Set-Text "E49" "-0.12%"

Set-Text "D50" "0.0001999"
Set-Text "E50" "-0.12%"
